$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39 (shifts the old row 39 down to row 40),
# duplicating row 38's formatting.
$ws.Rows.Item(39).Insert()

# New row 39 gets the values that row 38 used to have (before today's update).
$ws.Range("A39").Value = 1
$ws.Range("B39").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C39").Value = "Arica y Parinacota"
$ws.Range("D39").Value = 44172
$ws.Range("D39").NumberFormat = $ws.Range("D38").NumberFormat
$ws.Range("E39").Value = 15
$ws.Range("F39").Value = 100112044
$ws.Range("G39").Value = "Perejil"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 200
$ws.Range("K39").Value = 1300
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = 1400
$ws.Range("N39").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O39").Value = "Región de Arica y Parinacota"
$ws.Range("P39").Value = 700
$ws.Range("Q39").Value = 2
$ws.Range("R39").Value = "Hortaliza"

# Row 38 is updated in place with the new price data for this week.
$ws.Range("D38").Value = 45021
$ws.Range("J38").Value = 450
$ws.Range("K38").Value = 2000
$ws.Range("L38").Value = 2500
$ws.Range("M38").Value = 2111
$ws.Range("P38").Value = 1056
